$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D2:D51 to text format so numeric-looking price strings are preserved as text
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Apply cell value updates per the diff
$ws.Range("D2").Value = "69.167.20"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "3.886.60"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "529.45"
$ws.Range("E5").Value = "  +9.01%  "
$ws.Range("D6").Value = "144.08"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.717"
$ws.Range("E9").Value = "  -3.37%  "
$ws.Range("E10").Value = "  -5.51%  "
$ws.Range("E11").Value = "  -5.96%  "
$ws.Range("D12").Value = "41.90"
$ws.Range("E12").Value = "  -2.61%  "
$ws.Range("D13").Value = "4.512.04"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "10.20"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").Value = "3.887.71"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "13.98"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "1.22"
$ws.Range("E17").Value = "  +6.63%  "
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").Value = "20.29"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").Value = "69.165.51"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "423.17"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").Value = "3.38"
$ws.Range("E22").Value = "  -5.42%  "
$ws.Range("E23").Value = "  -4.08%  "
$ws.Range("D24").Value = "87.43"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("E25").Value = "  +8.23%  "
$ws.Range("D26").Value = "11.39"
$ws.Range("E26").Value = "  -7.64%  "
$ws.Range("D27").Value = "10.56"
$ws.Range("E27").Value = "  -3.65%  "
$ws.Range("D28").Value = "36.29"
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("D29").Value = "696.20"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("D30").Value = "13.18"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("D33").Value = "67.82"
$ws.Range("E33").Value = "  +10.02%  "
$ws.Range("D34").Value = "0.431"
$ws.Range("E34").Value = "  +6.40%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "5.93"
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0854"
$ws.Range("E36").Value = "  -5.01%  "
$ws.Range("D37").Value = "40.03"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "3.28"
$ws.Range("E41").Value = "  +5.80%  "
$ws.Range("D42").Value = "3.23"
$ws.Range("E42").Value = "  +6.66%  "
$ws.Range("D43").Value = "0.0480"
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("D44").Value = "2.78"
$ws.Range("E44").Value = "  -6.63%  "
$ws.Range("D45").Value = "3.39"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").Value = "2.98"
$ws.Range("E47").Value = "  +6.12%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.753.15"
$ws.Range("E48").Value = "  +15.01%  "
$ws.Range("D49").Value = "144.55"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "0.000268"
$ws.Range("E50").Value = "  +7.78%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "3.27"
$ws.Range("E51").Value = "  -2.99%  "

# Restore default (Normal) style for the price column so no stray number formats remain
$priceRange.Style = "Normal"

